$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.034.81"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "3.140.25"
$ws.Range("E3").Value = "  -0.50%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "602.57"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.57%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "143.12"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -3.27%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.133.95"
$ws.Range("E8").Value = "  -0.60%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.524"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("E10").Value = "  -1.53%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "5.39"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.54%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.467"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("E13").Value = "  -1.11%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "35.06"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").Value = "3.649.29"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D17").Value = "63.991.65"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "3.132.12"
$ws.Range("E18").Value = "  -0.70%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.83"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.35%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "487.51"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.50%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.71"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.709"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.06%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.65"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -4.28%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "86.98"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +3.48%  "
$ws.Range("E25").Value = "  -2.16%  "
$ws.Range("E26").Value = "  +0.05%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.76"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("E28").Value = "  -3.10%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "6.99"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("E31").Value = "  +2.54%  "
$ws.Range("E32").Value = "  -6.82%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("E35").Value = "  -3.08%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "6.01"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.01%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "52.53"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("E38").Value = "  -5.45%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.97"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -6.68%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "438.25"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -4.47%  "
$ws.Range("E41").Value = "  -1.46%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("D44").Value = "2.878.86"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("E45").Value = "  -3.12%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.21"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -5.18%  "
$ws.Range("E47").Value = "  -3.34%  "
$ws.Range("E48").Value = "  -0.06%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "25.98"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -2.01%  "
$ws.Range("E50").Value = "  -0.33%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "120.99"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
